$d = $word.ActiveDocument

$d.Content.Find.Execute("**Nonprofit**", $false, $false, $false, $false, $false, $true, 1, $false, "Shelter Care", 2)
